$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): set NumberFormat to Text first so numeric-looking
# strings (e.g. "1.70", "0.110") keep their literal digits instead of
# being auto-parsed into numbers; ClearFormats() afterwards restores the
# original (default) cell style while the stored value stays textual.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.534.10'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.239.60'
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.37'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.21'
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.01'
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.95'
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.22'
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.270.82'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.835'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.58'
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '44.284.72'
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0951'
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.85'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '65.49'
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.25'
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("D24").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.79'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.42'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.95'
$ws.Range("D29").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.43'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0796'
$ws.Range("D32").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.03'
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.110'
$ws.Range("D35").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.83'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.87'
$ws.Range("D38").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.808.72'
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.70'
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '78.77'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '70.58'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.77'
$ws.Range("D48").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.32'
$ws.Range("D51").ClearFormats()

# Volume(1h) column (E): plain text assignment is safe here (leading/
# trailing spaces and the percent sign prevent numeric coercion).
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("E26").Value = '  +4.36%  '
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("E28").Value = '  -1.25%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("E34").Value = '  -5.89%  '
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  +4.15%  '
$ws.Range("E38").Value = '  +0.44%  '
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("E43").Value = '  +5.26%  '
$ws.Range("E44").Value = '  +14.90%  '
$ws.Range("E45").Value = '  +3.32%  '
$ws.Range("E46").Value = '  -6.85%  '
$ws.Range("E47").Value = '  +3.39%  '
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("E50").Value = '  +0.17%  '
$ws.Range("E51").Value = '  +1.09%  '
